# [Kadastro App] Yeni kayit eklendi: 2925
# Appends the new record (Kayit No 2925) as a new last row on both the
# "Kayitlar" master sheet and the per-birim "Erdemli" sheet.

$wb = $excel.ActiveWorkbook

$recordValues = @("2925", "2025-09-08", "Erdemli", "1", "LİHKAB", "AYHAN KARADAYI (K.Teknisyeni), ÖZKAN AKBAŞ (Mühendis)")
$columns = @("A", "B", "C", "D", "E", "F")

function Add-KayitRow {
    param($Worksheet, $RowNumber)

    $rangeAddress = "A" + $RowNumber + ":F" + $RowNumber
    $rowRange = $Worksheet.Range($rangeAddress)
    # Every existing column in this workbook stores its data as plain text
    # (record numbers, parcel counts and dates are all text, not numbers/
    # dates) - force the Text format first so Excel doesn't silently
    # reinterpret "2925" / "2025-09-08" / "1" as numeric or date values.
    $rowRange.NumberFormat = "@"

    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellAddress = $columns[$i] + $RowNumber
        $Worksheet.Range($cellAddress).Value = $recordValues[$i]
    }
}

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
Add-KayitRow $wsKayitlar 18

$wsErdemli = $wb.Worksheets.Item("Erdemli")
Add-KayitRow $wsErdemli 17
